$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A32").Value = 29
$ws.Range("B32").Value = "2：57-5：15"
$ws.Range("C32").Value = "到了微分的定义"

$ws.Range("C32").Select()
